# Fixed a bug where Portugal was able to colonize new world too early if lied to.
# Update the "lied to" scenario table (rows 22-25) on sheet List1 with corrected
# base range / modifier values. Dependent formulas in columns E and H recalc
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("List1")

# Row 22 (Tenerife -> Verde, lied)
$ws.Range("B22").Value = 116
$ws.Range("D22").Value = 0.3
$ws.Range("G22").Value = 0.4

# Row 23 (Verde -> Africa, lied)
$ws.Range("B23").Value = 138
$ws.Range("D23").Value = 0.55000000000000004
$ws.Range("G23").Value = 0.72499999999999998

# Row 24 (Verde -> Brazil, lied)
$ws.Range("D24").Value = 0.6
$ws.Range("G24").Value = 0.75

# Row 25 (Range +50%, lied)
$ws.Range("D25").Value = 0.65
$ws.Range("G25").Value = 0.85

# Update the selected cell to reflect where the author ended up working
$ws.Range("E23").Select()

$excel.Calculate()
